# Build site at 2022-09-26 16:07:08 UTC
#
# The "Docentes responsaveis:" row (old row 13, which only had B/C values
# and no A-column label) is removed entirely, shifting every row below it
# up by one. A handful of the long descriptive paragraphs further down the
# sheet were then overwritten with short (and in a couple of cases clearly
# mismatched/reused) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 13 (Docentes responsaveis value row) entirely; this
# shifts rows 14-25 up to become rows 13-24.
$ws.Rows.Item(13).Delete()

# Row 10 (Objetivos:) - long paragraph replaced with the professor id/name.
$ws.Range("B10:C10").Value = "1643715 - Paulo Atsushi Suzuki"

# Row 13 (Programa resumido:, after the shift) - replaced with "Semestral".
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 (Programa:, after the shift) - replaced with the activation date.
# Use Copy from the existing "01/01/2012" text cell (row 8) instead of a
# plain Value assignment so Excel doesn't auto-convert the string into a
# date serial number - it must stay a literal text value.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# Row 18 (Metodo:, after the shift) - replaced with the professor id/name.
$ws.Range("B18:C18").Value = "1643715 - Paulo Atsushi Suzuki"

# Row 19 (Criterio:, after the shift) - now holds the "Aulas expositivas..." text.
$ws.Range("B19:C19").Value = "Aulas expositivas, seminários e exercícios comentados."

# Row 20 (Norma de recuperação:, after the shift) - now holds the "Média aritmética..." text.
$ws.Range("B20:C20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."

# Row 21 (Bibliografia:, after the shift) - now holds the "Aplicação de uma prova..." text.
$ws.Range("B21:C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
